# Prefix each worksheet's name onto the Step/command names in column A.
# Sheets that are "protocol" sheets (everything except the first 5 summary
# sheets: AshleyJourney, NRWaves, PersonalAshley, PositiveSpin, ReEngagement)
# get their own name + a space prepended to every existing value in column A,
# for all data rows (i.e. every row below the header row 1).

$wb = $excel.ActiveWorkbook

# Names of sheets that must NOT be touched by this edit.
$skipNames = @("AshleyJourney", "NRWaves", "PersonalAshley", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($skipNames -contains $ws.Name) {
        continue
    }

    $prefix = $ws.Name

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($null -ne $current -and "$current" -ne "") {
            $cell.Value = "$prefix " + "$current"
        }
    }
}
